# Update vm_pu.xlsx results for the 380 kV case: refresh the per-bus
# voltage magnitude (p.u.) values in res_bus/vm_pu.xlsx (rows 2-25,
# columns B-F and I-N). Column G (slack bus) stays at 1 p.u.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028483841947728
$ws.Range("D2").Value = 1.036293159674117
$ws.Range("E2").Value = 1.028430210264438
$ws.Range("F2").Value = 1.044930503471274
$ws.Range("I2").Value = 1.033052814160182
$ws.Range("J2").Value = 1.033635906682086
$ws.Range("K2").Value = 1.039087315974349
$ws.Range("L2").Value = 1.031247050133502
$ws.Range("M2").Value = 1.047700169207753
$ws.Range("N2").Value = 1.015169284544041
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029663487152199
$ws.Range("D3").Value = 1.037204709419734
$ws.Range("E3").Value = 1.029437944682024
$ws.Range("F3").Value = 1.046079220387699
$ws.Range("I3").Value = 1.033306826363643
$ws.Range("J3").Value = 1.034455000956907
$ws.Range("K3").Value = 1.039808166652206
$ws.Range("L3").Value = 1.032062179405888
$ws.Range("M3").Value = 1.048659344620985
$ws.Range("N3").Value = 1.015444788169815
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030426448082379
$ws.Range("D4").Value = 1.037793898047136
$ws.Range("E4").Value = 1.030090076115165
$ws.Range("F4").Value = 1.046822332074717
$ws.Range("I4").Value = 1.033469157221069
$ws.Range("J4").Value = 1.03498416772835
$ws.Range("K4").Value = 1.040273334609664
$ws.Range("L4").Value = 1.032589078856522
$ws.Range("M4").Value = 1.049279232802774
$ws.Range("N4").Value = 1.015622654701629
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030747115131232
$ws.Range("D5").Value = 1.038041438964668
$ws.Range("E5").Value = 1.03036424718084
$ws.Range("F5").Value = 1.047134692866429
$ws.Range("I5").Value = 1.033536914721989
$ws.Range("J5").Value = 1.035206428633012
$ws.Range("K5").Value = 1.040468587517239
$ws.Range("L5").Value = 1.032810457892201
$ws.Range("M5").Value = 1.049539652010974
$ws.Range("N5").Value = 1.015697333564818
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030800951796038
$ws.Range("D6").Value = 1.038082993156886
$ws.Range("E6").Value = 1.030410282608972
$ws.Range("F6").Value = 1.047187137101683
$ws.Range("I6").Value = 1.033548262984467
$ws.Range("J6").Value = 1.03524373546377
$ws.Range("K6").Value = 1.040501353544441
$ws.Range("L6").Value = 1.032847620832164
$ws.Range("M6").Value = 1.049583366903653
$ws.Range("N6").Value = 1.015709866837221
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030430733170218
$ws.Range("D7").Value = 1.037797206305514
$ws.Range("E7").Value = 1.03009373954255
$ws.Range("F7").Value = 1.046826506024155
$ws.Range("I7").Value = 1.033470064511219
$ws.Range("J7").Value = 1.034987138377524
$ws.Range("K7").Value = 1.040275944781325
$ws.Range("L7").Value = 1.032592037442033
$ws.Range("M7").Value = 1.049282713250353
$ws.Range("N7").Value = 1.015623652942128
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02888258171915
$ws.Range("D8").Value = 1.036601355997819
$ws.Range("E8").Value = 1.028770766839494
$ws.Range("F8").Value = 1.045318756587515
$ws.Range("I8").Value = 1.033139079567337
$ws.Range("J8").Value = 1.03391289823332
$ws.Range("K8").Value = 1.039331194226883
$ws.Range("L8").Value = 1.031522639867315
$ws.Range("M8").Value = 1.048024485352548
$ws.Range("N8").Value = 1.015262475788121
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026151796882035
$ws.Range("D9").Value = 1.034489154388153
$ws.Range("E9").Value = 1.026439947457747
$ws.Range("F9").Value = 1.042660450603975
$ws.Range("I9").Value = 1.032540274692272
$ws.Range("J9").Value = 1.032013465955649
$ws.Range("K9").Value = 1.037656668628225
$ws.Range("L9").Value = 1.0296340300203
$ws.Range("M9").Value = 1.045801453127023
$ws.Range("N9").Value = 1.01462294226343
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024329319693108
$ws.Range("D10").Value = 1.033077652422169
$ws.Range("E10").Value = 1.024886308782953
$ws.Range("F10").Value = 1.040887203174093
$ws.Range("I10").Value = 1.032130596914006
$ws.Range("J10").Value = 1.030742766288869
$ws.Range("K10").Value = 1.036533724796439
$ws.Range("L10").Value = 1.028372090944785
$ws.Range("M10").Value = 1.044315430744241
$ws.Range("N10").Value = 1.014194493556303
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023539678368352
$ws.Range("D11").Value = 1.032465649436868
$ws.Range("E11").Value = 1.024213611900726
$ws.Range("F11").Value = 1.04011910492073
$ws.Range("I11").Value = 1.031950715019117
$ws.Range("J11").Value = 1.030191480485536
$ws.Range("K11").Value = 1.036045905061169
$ws.Range("L11").Value = 1.02782496666098
$ws.Range("M11").Value = 1.043671005138318
$ws.Range("N11").Value = 1.01400847080193
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023246293527671
$ws.Range("D12").Value = 1.03223820124517
$ws.Range("E12").Value = 1.023963747091177
$ws.Range("F12").Value = 1.039833756980762
$ws.Range("I12").Value = 1.031883524668178
$ws.Range("J12").Value = 1.029986547287091
$ws.Range("K12").Value = 1.035864469193859
$ws.Range("L12").Value = 1.027621634677396
$ws.Range("M12").Value = 1.043431490122668
$ws.Range("N12").Value = 1.013939297977884
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023309229088149
$ws.Range("D13").Value = 1.03228699523555
$ws.Range("E13").Value = 1.024017343742897
$ws.Range("F13").Value = 1.0398949670083
$ws.Range("I13").Value = 1.031897954165119
$ws.Range("J13").Value = 1.030030513463241
$ws.Range("K13").Value = 1.035903398582807
$ws.Range("L13").Value = 1.027665254876521
$ws.Range("M13").Value = 1.043482873562393
$ws.Range("N13").Value = 1.01395413921386
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.0235154286677
$ws.Range("D14").Value = 1.032446851013148
$ws.Range("E14").Value = 1.024192957906726
$ws.Range("F14").Value = 1.040095518837057
$ws.Range("I14").Value = 1.031945168679226
$ws.Range("J14").Value = 1.030174543936787
$ws.Range("K14").Value = 1.036030912373909
$ws.Range("L14").Value = 1.027808161337386
$ws.Range("M14").Value = 1.043651209751088
$ws.Range("N14").Value = 1.014002754500669
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02364246477673
$ws.Range("D15").Value = 1.032545327105116
$ws.Range("E15").Value = 1.024301160187844
$ws.Range("F15").Value = 1.040219079823135
$ws.Range("I15").Value = 1.031974209504479
$ws.Range("J15").Value = 1.030263264481649
$ws.Range("K15").Value = 1.036109446279041
$ws.Range("L15").Value = 1.027896196691475
$ws.Range("M15").Value = 1.0437549077636
$ws.Range("N15").Value = 1.014032697943887
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024381714928959
$ws.Range("D16").Value = 1.033118251837473
$ws.Range("E16").Value = 1.024930954187304
$ws.Range("F16").Value = 1.040938173507044
$ws.Range("I16").Value = 1.032142482622759
$ws.Range("J16").Value = 1.03077933078438
$ws.Range("K16").Value = 1.03656606647214
$ws.Range("L16").Value = 1.028408387021396
$ws.Range("M16").Value = 1.04435817864981
$ws.Range("N16").Value = 1.014206828677267
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024845292441718
$ws.Range("D17").Value = 1.033477414046184
$ws.Range("E17").Value = 1.025326016981519
$ws.Range("F17").Value = 1.041389168578105
$ws.Range("I17").Value = 1.032247369391759
$ws.Range("J17").Value = 1.031102759723314
$ws.Range("K17").Value = 1.036852069260923
$ws.Range("L17").Value = 1.028729483317361
$ws.Range("M17").Value = 1.04473633450075
$ws.Range("N17").Value = 1.014315921754683
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025115641595198
$ws.Range("D18").Value = 1.033686828871238
$ws.Range("E18").Value = 1.025556454255176
$ws.Range("F18").Value = 1.041652200379871
$ws.Range("I18").Value = 1.032308307966475
$ws.Range("J18").Value = 1.031291307727646
$ws.Range("K18").Value = 1.03701873768962
$ws.Range("L18").Value = 1.028916706240148
$ws.Range("M18").Value = 1.044956813085657
$ws.Range("N18").Value = 1.014379505522879
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025207815636188
$ws.Range("D19").Value = 1.033758220639089
$ws.Range("E19").Value = 1.025635028092171
$ws.Range("F19").Value = 1.041741883060364
$ws.Range("I19").Value = 1.032329045728867
$ws.Range("J19").Value = 1.031355580353016
$ws.Range("K19").Value = 1.037075541533715
$ws.Range("L19").Value = 1.02898053306596
$ws.Range("M19").Value = 1.045031974791673
$ws.Range("N19").Value = 1.014401177748542
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024795559920184
$ws.Range("D20").Value = 1.033438887458016
$ws.Range("E20").Value = 1.025283630095133
$ws.Range("F20").Value = 1.041340783775071
$ws.Range("I20").Value = 1.032236140871683
$ws.Range("J20").Value = 1.031068069483301
$ws.Range("K20").Value = 1.036821399614017
$ws.Range("L20").Value = 1.028695039673923
$ws.Range("M20").Value = 1.044695771619479
$ws.Range("N20").Value = 1.014304222106537
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023454710171873
$ws.Range("D21").Value = 1.032399780874093
$ws.Range("E21").Value = 1.024141243789769
$ws.Range("F21").Value = 1.040036462494355
$ws.Range("I21").Value = 1.031931275514014
$ws.Range("J21").Value = 1.03013213500902
$ws.Range("K21").Value = 1.035993369314836
$ws.Range("L21").Value = 1.027766081871805
$ws.Range("M21").Value = 1.043601642997253
$ws.Range("N21").Value = 1.013988440602461
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022611218699496
$ws.Range("D22").Value = 1.031745741063157
$ws.Range("E22").Value = 1.023423007200635
$ws.Range("F22").Value = 1.039216139940756
$ws.Range("I22").Value = 1.031737429367099
$ws.Range("J22").Value = 1.029542743161602
$ws.Range("K22").Value = 1.035471376959309
$ws.Range("L22").Value = 1.027181397368105
$ws.Range("M22").Value = 1.042912871814952
$ws.Range("N22").Value = 1.013789458152906
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023058412527923
$ws.Range("D23").Value = 1.032092527740149
$ws.Range("E23").Value = 1.023803755808717
$ws.Range("F23").Value = 1.039651032053776
$ws.Range("I23").Value = 1.031840396216724
$ws.Range("J23").Value = 1.02985527973768
$ws.Range("K23").Value = 1.035748225787004
$ws.Range("L23").Value = 1.027491408035784
$ws.Range("M23").Value = 1.043278083340476
$ws.Range("N23").Value = 1.013894984149096
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024818032053579
$ws.Range("D24").Value = 1.03345629620728
$ws.Range("E24").Value = 1.025302782892861
$ws.Range("F24").Value = 1.0413626468653
$ws.Range("I24").Value = 1.032241215298523
$ws.Range("J24").Value = 1.03108374482671
$ws.Range("K24").Value = 1.036835258377539
$ws.Range("L24").Value = 1.028710603481245
$ws.Range("M24").Value = 1.044714100527986
$ws.Range("N24").Value = 1.014309508823449
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026858106600247
$ws.Range("D25").Value = 1.035035800270202
$ws.Range("E25").Value = 1.027042474228031
$ws.Range("F25").Value = 1.043347865700439
$ws.Range("I25").Value = 1.032696925690156
$ws.Range("J25").Value = 1.032505288429396
$ws.Range("K25").Value = 1.038090733154679
$ws.Range("L25").Value = 1.030122782687475
$ws.Range("M25").Value = 1.046376861601671
$ws.Range("N25").Value = 1.014788644931483
